$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows (2-9) with corrected/re-simulated values ---

# Row 2 (Q0)
$ws.Range("B2").Value = 0.1196643795433455
$ws.Range("C2").Value = 0.4311790077471714
$ws.Range("D2").Value = 0.2593780065637157
$ws.Range("E2").Value = 0.509291671406195
$ws.Range("F2").Value = 0.5137208310754531

# Row 3 (Q1)
$ws.Range("B3").Value = 0.1031233158734826
$ws.Range("C3").Value = 0.3615473318654818
$ws.Range("D3").Value = 0.2153226717276779
$ws.Range("E3").Value = 0.4640287401957749
$ws.Range("F3").Value = 0.4708987236181263

# Row 4 (Q2)
$ws.Range("B4").Value = 0.1715531831306264
$ws.Range("C4").Value = 0.6362730885198975
$ws.Range("D4").Value = 0.7952461905435021
$ws.Range("E4").Value = 0.8917657711212638
$ws.Range("F4").Value = 0.914021501180113

# Row 5 (Q3)
$ws.Range("B5").Value = 0.2524150092910198
$ws.Range("C5").Value = 0.6280498924479758
$ws.Range("D5").Value = 1.020506542609963
$ws.Range("E5").Value = 1.010201238669783
$ws.Range("F5").Value = 1.025900836467168

# Row 6 (Q4)
$ws.Range("B6").Value = 0.1672632190142579
$ws.Range("C6").Value = 0.4920935882630871
$ws.Range("D6").Value = 0.5780268184013607
$ws.Range("E6").Value = 0.7602807497243111
$ws.Range("F6").Value = 0.7817713746261321

# Row 7 (Q5)
$ws.Range("B7").Value = 0.04919676379112856
$ws.Range("C7").Value = 0.4710688180516911
$ws.Range("D7").Value = 0.4742404420980256
$ws.Range("E7").Value = 0.6886511759214716
$ws.Range("F7").Value = 0.7285586013470837
$ws.Range("G7").Value = 9

# Row 8 (Q6)
$ws.Range("B8").Value = -0.09646499961979775
$ws.Range("C8").Value = 0.4357283444034626
$ws.Range("D8").Value = 0.2943364263897597
$ws.Range("E8").Value = 0.5425278116279014
$ws.Range("F8").Value = 0.584839393582319
$ws.Range("G8").Value = 6

# Row 9 (Q7)
$ws.Range("B9").Value = -0.1901586249338153
$ws.Range("C9").Value = 0.3356701559601163
$ws.Range("D9").Value = 0.2161827050533335
$ws.Range("E9").Value = 0.4649545193385408
$ws.Range("F9").Value = 0.5196475763677928
$ws.Range("G9").Value = 3

# --- Append new row 10 (Q8) ---
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A10").Value = "Q8"
$ws.Range("B10").Value = -0.4377465499738378
$ws.Range("C10").Value = 0.4377465499738378
$ws.Range("D10").Value = 0.1916220420139977
$ws.Range("E10").Value = 0.4377465499738378
$ws.Range("G10").Value = 1
